# Applies the "update from upstream" translation refresh for the
# Combat Extended Armors sheet:
#   - Korean label/description for the "large backpack" item is updated
#     from "백팩" (backpack) wording to "배낭" (knapsack) wording.
#   - A conditional format is added on F2:F3 that highlights the Korean
#     translation cell green whenever it still equals the English source
#     string (i.e. still untranslated).
#   - The sheet is renamed to reflect the newer export date, and the
#     last active-cell selection is moved to match upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Korean translation text updates ---------------------------------
$ws.Range("F3").Value = "등산용 대형 배낭입니다. 산업계 군대의 보병에게 지급되는 종류의 것과 기능상으로 동일합니다."
$ws.Range("F2").Value = "대형 배낭"

# --- Conditional formatting: flag untranslated cells in green --------
$rng = $ws.Range("F2:F3")
$fc = $rng.FormatConditions.Add(2, 3, "(E2=F2)")
$fc.Interior.Color = 5296274

# --- Sheet rename (newer upstream export date) ------------------------
$ws.Name = "Main_240603"

# --- Restore the last saved selection ---------------------------------
$ws.Range("F25").Select()
